# Apply updated "想去人数" (interested count) figures across the workbook.
# Source data changed upstream; F-column values below are bumped by the
# same amounts on the sheets that aggregate/repeat the same events.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 8020
$ws1.Range("F4").Value = 1890
$ws1.Range("F9").Value = 39
$ws1.Range("F16").Value = 154
$ws1.Range("F25").Value = 1
$ws1.Range("F28").Value = 9
$ws1.Range("F30").Value = 839

# --- Sheet "本地生活" (local life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 2311

# --- Sheet "全部类型" (all categories, aggregates the above events) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 2311
$ws4.Range("F5").Value = 8020
$ws4.Range("F7").Value = 1890
$ws4.Range("F13").Value = 39
$ws4.Range("F32").Value = 839
